$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.399.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.711.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5349"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06617"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.97"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.560"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.711.53"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.949.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5773"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8194"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.93"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.386.85"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.18"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.671"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.985"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.727"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.285"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05406"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.290"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.497"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.643"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.878"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9508"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5862"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.864"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.045.34"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8404"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.92"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.854.73"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.04"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.089"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05245"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.68%  "
